$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.770.75"
$ws.Range("E2").Value = "  +1.26%  "
$ws.Range("D3").Value = "2.699.34"
$ws.Range("E3").Value = "  +2.66%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "608.81"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.26%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "157.79"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.77%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.589"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.20%  "
$ws.Range("E9").Value = "  +5.41%  "
$ws.Range("E10").Value = "  +3.91%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.403"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.72%  "
$ws.Range("E13").Value = "  +4.64%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000204"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +10.29%  "
$ws.Range("D15").Value = "3.185.87"
$ws.Range("E15").Value = "  +2.72%  "
$ws.Range("D16").Value = "65.654.45"
$ws.Range("E16").Value = "  +1.18%  "
$ws.Range("D17").Value = "2.708.58"
$ws.Range("E17").Value = "  +2.53%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.72"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.79%  "
$ws.Range("E19").Value = "  +2.32%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "358.98"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.19%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.64"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.02%  "
$ws.Range("E22").Value = "  -0.12%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.97"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.96%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.94"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.86%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000106"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +12.33%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.66"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.65%  "
$ws.Range("E27").Value = "  +3.43%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.53"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +5.65%  "
$ws.Range("E29").Value = "  +4.18%  "
$ws.Range("E30").Value = "  +5.41%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "545.42"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.11%  "
$ws.Range("E32").Value = "  +0.05%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.82"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.20%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.71"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +6.82%  "
$ws.Range("E35").Value = "  -1.64%  "
$ws.Range("E36").Value = "  +2.33%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "20.88"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.21%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "163.30"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.18%  "
$ws.Range("E39").Value = "  +0.46%  "
$ws.Range("E40").Value = "  -0.04%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "172.33"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.58%  "
$ws.Range("E42").Value = "  -0.05%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "42.60"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.75%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.20"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.01%  "
$ws.Range("E45").Value = "  +0.79%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "23.69"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.06%  "
$ws.Range("E47").Value = "  +3.66%  "
$ws.Range("E48").Value = "  +5.10%  "
$ws.Range("E49").Value = "  +1.53%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "21.08"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +8.61%  "
$ws.Range("E51").Value = "  +1.51%  "
